$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates ---
# These values look like plain numbers (e.g. "21.32") to Excel's parser,
# so we force the cell to Text format before assigning the value, then
# restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.686.52'
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.532.51'
$ws.Range("D3").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '205.53'
$ws.Range("D5").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '21.32'
$ws.Range("D8").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0578'
$ws.Range("D10").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.752.01'
$ws.Range("D12").Style = "Normal"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.546.80'
$ws.Range("D13").Style = "Normal"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.505'
$ws.Range("D15").Style = "Normal"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.694.09'
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.19'
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '212.03'
$ws.Range("D18").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.19'
$ws.Range("D20").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.06'
$ws.Range("D23").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.75'
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.56'
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.81'
$ws.Range("D27").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0452'
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.23'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.356.82'
$ws.Range("D33").Style = "Normal"

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.49'
$ws.Range("D35").Style = "Normal"

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.934'
$ws.Range("D37").Style = "Normal"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.796'
$ws.Range("D41").Style = "Normal"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '62.37'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.665.42'
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.27'
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0504'
$ws.Range("D49").Style = "Normal"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₇0970'
$ws.Range("D50").Style = "Normal"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0942'
$ws.Range("D51").Style = "Normal"

# --- Column E (Volume(1h)) updates ---
# These strings contain surrounding spaces and a percent sign, so Excel
# keeps them as plain text automatically.
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("E3").Value = '  -1.68%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -3.03%  '
$ws.Range("E9").Value = '  -1.31%  '
$ws.Range("E10").Value = '  -0.62%  '
$ws.Range("E11").Value = '  -1.12%  '
$ws.Range("E12").Value = '  -1.48%  '
$ws.Range("E13").Value = '  -0.75%  '
$ws.Range("E14").Value = '  -1.86%  '
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("E19").Value = '  +0.92%  '
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("E22").Value = '  -2.14%  '
$ws.Range("E23").Value = '  -3.14%  '
$ws.Range("E24").Value = '  -2.40%  '
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("E26").Value = '  -3.11%  '
$ws.Range("E27").Value = '  -0.10%  '
$ws.Range("E28").Value = '  -0.05%  '
$ws.Range("E29").Value = '  -1.33%  '
$ws.Range("E30").Value = '  -0.87%  '
$ws.Range("E31").Value = '  -2.24%  '
$ws.Range("E32").Value = '  +2.50%  '
$ws.Range("E33").Value = '  -2.09%  '
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("E35").Value = '  -3.68%  '
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("E39").Value = '  +0.74%  '
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("E42").Value = '  +5.23%  '
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("E45").Value = '  -1.55%  '
$ws.Range("E46").Value = '  -1.18%  '
$ws.Range("E47").Value = '  -1.56%  '
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("E49").Value = '  +2.29%  '
$ws.Range("E50").Value = '  -0.38%  '
$ws.Range("E51").Value = '  -0.61%  '
